$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "68.368.53"
Set-TextValue "E2" "  +0.87%  "
Set-TextValue "D3" "2.453.20"
Set-TextValue "E3" "  +1.02%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "557.59"
Set-TextValue "E5" "  +0.96%  "
Set-TextValue "D6" "162.59"
Set-TextValue "E6" "  +1.70%  "
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "D8" "0.510"
Set-TextValue "E8" "  +2.76%  "
Set-TextValue "D9" "2.451.86"
Set-TextValue "E9" "  +0.99%  "
Set-TextValue "E10" "  +7.45%  "
Set-TextValue "E11" "  +0.63%  "
Set-TextValue "E12" "  +2.40%  "
Set-TextValue "D13" "0.328"
Set-TextValue "E13" "  -1.08%  "
Set-TextValue "D14" "68.275.61"
Set-TextValue "E14" "  +0.73%  "
Set-TextValue "E15" "  +3.11%  "
Set-TextValue "E16" "  +2.53%  "
Set-TextValue "D17" "10.48"
Set-TextValue "E17" "  -1.77%  "
Set-TextValue "D18" "336.57"
Set-TextValue "E18" "  +0.04%  "
Set-TextValue "D19" "6.88"
Set-TextValue "E19" "  -1.41%  "
Set-TextValue "E20" "  +2.54%  "
Set-TextValue "E21" "  +0.05%  "
Set-TextValue "E22" "  +3.23%  "
Set-TextValue "D23" "66.61"
Set-TextValue "E23" "  +1.26%  "
Set-TextValue "D24" "3.66"
Set-TextValue "E24" "  +2.09%  "
Set-TextValue "D25" "8.11"
Set-TextValue "E25" "  +2.13%  "
Set-TextValue "D26" "0.0₃0814"
Set-TextValue "E26" "  +0.94%  "
Set-TextValue "D27" "7.19"
Set-TextValue "E27" "  +2.87%  "
Set-TextValue "D28" "1.00"
Set-TextValue "E28" "  +0.05%  "
Set-TextValue "D29" "423.76"
Set-TextValue "E29" "  +3.05%  "
Set-TextValue "E30" "  +3.45%  "
Set-TextValue "E31" "  +0.45%  "
Set-TextValue "E32" "  +2.37%  "
Set-TextValue "D33" "18.98"
Set-TextValue "E33" "  +0.07%  "
Set-TextValue "E34" "  +0.02%  "
Set-TextValue "D35" "17.77"
Set-TextValue "E35" "  +1.11%  "
Set-TextValue "E36" "  -1.03%  "
Set-TextValue "E37" "  +2.92%  "
Set-TextValue "E38" "  -0.75%  "
Set-TextValue "E39" "  +1.31%  "
Set-TextValue "E40" "  +0.41%  "
Set-TextValue "E41" "  +0.52%  "
Set-TextValue "E42" "  +2.69%  "
Set-TextValue "D43" "129.55"
Set-TextValue "E43" "  -1.81%  "
Set-TextValue "D44" "0.0719"
Set-TextValue "E44" "  +1.46%  "
Set-TextValue "E45" "  +2.96%  "
Set-TextValue "E46" "  +1.98%  "
Set-TextValue "E47" "  +2.10%  "
Set-TextValue "E48" "  +1.35%  "
Set-TextValue "E49" "  -0.27%  "
Set-TextValue "B50" "THORChain"
Set-TextValue "C50" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D50" "4.88"
Set-TextValue "E50" "  -3.77%  "
Set-TextValue "B51" "InjectiveProtocol"
Set-TextValue "C51" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D51" "16.68"
Set-TextValue "E51" "  +2.19%  "
